# Insert a new data row above the current row 76, shifting the existing
# rows 76-103 down to 77-104 (dimension grows from R103 to R104), then
# populate the newly inserted row 76 with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Cells.Item(76, 1).Value  = 5
$ws.Cells.Item(76, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(76, 3).Value  = "Maule"
$ws.Cells.Item(76, 4).Value  = 44855
$ws.Cells.Item(76, 5).Value  = 7
$ws.Cells.Item(76, 6).Value  = 100112022
$ws.Cells.Item(76, 7).Value  = "Arveja Verde"
$ws.Cells.Item(76, 8).Value  = "Sin especificar"
$ws.Cells.Item(76, 9).Value  = "Primera"
$ws.Cells.Item(76, 10).Value = 150
$ws.Cells.Item(76, 11).Value = 13000
$ws.Cells.Item(76, 12).Value = 13000
$ws.Cells.Item(76, 13).Value = 13000
$ws.Cells.Item(76, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(76, 15).Value = "Región del Maule"
$ws.Cells.Item(76, 16).Value = 520
$ws.Cells.Item(76, 17).Value = 25
$ws.Cells.Item(76, 18).Value = "Hortaliza"
